$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" rows (16-22) are being re-sorted from descending
# (2210..2204) to ascending (2204..2210) order, the "Valor Mora" (F) value
# stays attached to its period (2210 -> 32707, all others -> 40000), and the
# "Salario Basico" (G) is updated from 1000000 to 908526 for every period.

$periods = @(2204, 2205, 2206, 2207, 2208, 2209, 2210)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value = $period

    if ($period -eq 2210) {
        $ws.Cells.Item($row, 6).Value = 32707
    } else {
        $ws.Cells.Item($row, 6).Value = 40000
    }

    $ws.Cells.Item($row, 7).Value = 908526
}
